# The workbook originally has three sheets:
#   strategy_id-0, strategy_id-6000, strategy_id-6002
# The target state has two sheets:
#   strategy_id-0, strategy_id-6081
# i.e. sheet "strategy_id-6002" is removed entirely, and
# sheet "strategy_id-6000" is renamed to "strategy_id-6081".

$wb = $excel.ActiveWorkbook

# Remove the obsolete worksheet (suppress the boolean return value).
$wb.Sheets.Item("strategy_id-6002").Delete() | Out-Null

# Rename the remaining strategy sheet to its new id.
$wb.Sheets.Item("strategy_id-6000").Name = "strategy_id-6081"
